# Update the per-seed / per-patient accuracy figures in the balanced-subjects
# CSETNet results sheet (visualization-of-results refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.91
$ws.Range("C2").Value = 0.99
$ws.Range("E2").Value = 0.71
$ws.Range("F2").Value = 0.9399999999999999
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 0.84
$ws.Range("I2").Value = 0.74
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 0.9
$ws.Range("B3").Value = 0.9399999999999999
$ws.Range("E3").Value = 0.51
$ws.Range("F3").Value = 0.95
$ws.Range("H3").Value = 0.82
$ws.Range("I3").Value = 0.8100000000000001
$ws.Range("K3").Value = 0.89
$ws.Range("B4").Value = 0.84
$ws.Range("C4").Value = 0.71
$ws.Range("E4").Value = 0.46
$ws.Range("F4").Value = 0.92
$ws.Range("H4").Value = 0.84
$ws.Range("I4").Value = 0.77
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 0.84
$ws.Range("E5").Value = 0.77
$ws.Range("F5").Value = 0.91
$ws.Range("H5").Value = 0.82
$ws.Range("I5").Value = 0.78
$ws.Range("K5").Value = 0.92
$ws.Range("B6").Value = 0.98
$ws.Range("E6").Value = 0.71
$ws.Range("F6").Value = 0.83
$ws.Range("H6").Value = 0.92
$ws.Range("I6").Value = 0.66
$ws.Range("J6").Value = 0.99
$ws.Range("K6").Value = 0.9
$ws.Range("B7").Value = 0.97
$ws.Range("C7").Value = 0.43
$ws.Range("E7").Value = 0.86
$ws.Range("F7").Value = 0.89
$ws.Range("H7").Value = 0.95
$ws.Range("I7").Value = 0.8
$ws.Range("B8").Value = 0.99
$ws.Range("C8").Value = 1
$ws.Range("E8").Value = 0.92
$ws.Range("F8").Value = 0.9
$ws.Range("G8").Value = 0.99
$ws.Range("H8").Value = 0.9
$ws.Range("I8").Value = 0.79
$ws.Range("J8").Value = 0.99
$ws.Range("K8").Value = 0.9399999999999999
$ws.Range("C9").Value = 0.97
$ws.Range("E9").Value = 0.84
$ws.Range("F9").Value = 0.93
$ws.Range("G9").Value = 0.99
$ws.Range("H9").Value = 0.85
$ws.Range("I9").Value = 0.68
$ws.Range("J9").Value = 0.99
$ws.Range("K9").Value = 0.91
$ws.Range("B10").Value = 0.91
$ws.Range("E10").Value = 0.52
$ws.Range("F10").Value = 0.91
$ws.Range("G10").Value = 0.99
$ws.Range("H10").Value = 0.93
$ws.Range("I10").Value = 0.76
$ws.Range("J10").Value = 0.99
$ws.Range("K10").Value = 0.89
$ws.Range("E11").Value = 0.92
$ws.Range("F11").Value = 0.9399999999999999
$ws.Range("H11").Value = 0.78
$ws.Range("I11").Value = 0.78
$ws.Range("J11").Value = 0.99
$ws.Range("B12").Value = 0.96
$ws.Range("C12").Value = 0.99
$ws.Range("E12").Value = 0.59
$ws.Range("F12").Value = 0.75
$ws.Range("G12").Value = 1
$ws.Range("H12").Value = 0.85
$ws.Range("I12").Value = 0.8
$ws.Range("K12").Value = 0.88
$ws.Range("B13").Value = 0.92
$ws.Range("C13").Value = 0.88
$ws.Range("E13").Value = 0.83
$ws.Range("F13").Value = 0.97
$ws.Range("H13").Value = 0.73
$ws.Range("I13").Value = 0.73
$ws.Range("J13").Value = 1
$ws.Range("K13").Value = 0.89
$ws.Range("B14").Value = 0.96
$ws.Range("C14").Value = 0.97
$ws.Range("E14").Value = 0.66
$ws.Range("F14").Value = 0.95
$ws.Range("G14").Value = 0.96
$ws.Range("H14").Value = 0.82
$ws.Range("I14").Value = 0.73
$ws.Range("J14").Value = 1
$ws.Range("K14").Value = 0.89
$ws.Range("B15").Value = 0.97
$ws.Range("C15").Value = 0.99
$ws.Range("E15").Value = 0.8
$ws.Range("F15").Value = 0.85
$ws.Range("G15").Value = 0.98
$ws.Range("H15").Value = 0.93
$ws.Range("I15").Value = 0.79
$ws.Range("K15").Value = 0.92
$ws.Range("C16").Value = 0.92
$ws.Range("E16").Value = 0.72
$ws.Range("F16").Value = 0.9
$ws.Range("I16").Value = 0.76
$ws.Range("J16").Value = 1
$ws.Range("K16").Value = 0.9
